$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting existing rows 2-3 down to 3-4
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with part numbers (columns L through AD)
$ws.Range("L2").Value = "E2-0052"
$ws.Range("M2").Value = "E1-0199"
$ws.Range("N2").Value = "E2-0047"
$ws.Range("O2").Value = "PC-1220"
$ws.Range("P2").Value = "PM-1006-SS"
$ws.Range("Q2").Value = "UA-1212"
$ws.Range("R2").Value = "E1-2530"
$ws.Range("S2").Value = "E2-0166"
$ws.Range("T2").Value = "E2-0177"
$ws.Range("U2").Value = "E2-0545"
$ws.Range("V2").Value = "E2-0154"
$ws.Range("W2").Value = "E2-0611"
$ws.Range("X2").Value = "BE9-2513"
$ws.Range("Y2").Value = "BE9-2513"
$ws.Range("Z2").Value = "BE9-2512"
$ws.Range("AA2").Value = "BE9-2511"
$ws.Range("AB2").Value = "BE9-2515"
$ws.Range("AC2").Value = "BE9-2514"
$ws.Range("AD2").Value = "BE9-2578"
